$d = $word.ActiveDocument

$target  = "Updating to version 2"
$newText = "Updating to version 3 from devops_test"

# Locate the paragraph that holds the marker sentence.
$rng = $d.Content
$found = $rng.Find.Execute($target, $true, $true, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text '$target'"
}

$para = $rng.Paragraphs(1)

# Split the paragraph in two so the new sentence lands in its own
# paragraph, right where the diff puts it (just before the trailing
# bookmark that closes out the original paragraph). Splitting flush
# against the very end of the match (i.e. right on the paragraph mark)
# pushes the new paragraph past any trailing bookmarks instead of
# landing before them, so split one character earlier and repair the
# text of both resulting paragraphs afterwards.
$splitAt = $rng.End - 1
$d.Range($splitAt, $splitAt).InsertParagraphAfter() | Out-Null

# Restore the first (original) paragraph's text.
$firstRange = $para.Range
$firstRange.MoveEnd(1, -1) | Out-Null
$firstRange.Text = $target

# Fill in the newly created paragraph (it inherited the trailing
# bookmark(s) that used to close out the original paragraph).
$secondPara = $para.Next()
$secondRange = $secondPara.Range
$secondRange.MoveEnd(1, -1) | Out-Null
$secondRange.Text = $newText
